$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) contain numeric-looking text that must remain exact text
# (preserving trailing zeros / multi-dot European-style formatting),
# so force text number format before assigning values.
$dCells = @("D2","D3","D5","D6","D8","D10","D13","D14","D15","D17","D18","D20","D21","D22","D24","D30","D32","D33","D34","D39","D40","D43","D44","D48")
foreach ($addr in $dCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range('D2').Value = '64.085.17'
$ws.Range('E2').Value = '  -0.24%  '
$ws.Range('D3').Value = '2.757.33'
$ws.Range('E3').Value = '  -1.04%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '575.71'
$ws.Range('E5').Value = '  -2.25%  '
$ws.Range('D6').Value = '159.19'
$ws.Range('E6').Value = '  -1.41%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').Value = '0.600'
$ws.Range('E8').Value = '  -3.25%  '
$ws.Range('E9').Value = '  -4.16%  '
$ws.Range('D10').Value = '5.90'
$ws.Range('E10').Value = '  -13.33%  '
$ws.Range('E11').Value = '  +3.36%  '
$ws.Range('E12').Value = '  -3.60%  '
$ws.Range('D13').Value = '3.247.03'
$ws.Range('E13').Value = '  -1.00%  '
$ws.Range('D14').Value = '26.95'
$ws.Range('E14').Value = '  -2.63%  '
$ws.Range('D15').Value = '63.712.93'
$ws.Range('E15').Value = '  -0.67%  '
$ws.Range('E16').Value = '  -5.57%  '
$ws.Range('D17').Value = '2.763.03'
$ws.Range('E17').Value = '  -0.82%  '
$ws.Range('D18').Value = '12.16'
$ws.Range('E18').Value = '  -2.17%  '
$ws.Range('E19').Value = '  -5.25%  '
$ws.Range('D20').Value = '359.16'
$ws.Range('E20').Value = '  -2.39%  '
$ws.Range('D21').Value = '6.65'
$ws.Range('E21').Value = '  -6.40%  '
$ws.Range('D22').Value = '0.998'
$ws.Range('E22').Value = '  -0.65%  '
$ws.Range('E23').Value = '  -8.62%  '
$ws.Range('D24').Value = '65.06'
$ws.Range('E25').Value = '  -3.90%  '
$ws.Range('E26').Value = '  -3.95%  '
$ws.Range('E27').Value = '  -0.09%  '
$ws.Range('E28').Value = '  -7.09%  '
$ws.Range('E29').Value = '  +1.16%  '
$ws.Range('D30').Value = '1.36'
$ws.Range('E30').Value = '  +7.03%  '
$ws.Range('E31').Value = '  -4.50%  '
$ws.Range('D32').Value = '170.15'
$ws.Range('E32').Value = '  -1.18%  '
$ws.Range('B33').Value = 'NEARProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D33').Value = '4.94'
$ws.Range('E33').Value = '  -4.71%  '
$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').Value = '20.23'
$ws.Range('E34').Value = '  -3.31%  '
$ws.Range('E35').Value = '  -1.90%  '
$ws.Range('E36').Value = '  +0.04%  '
$ws.Range('E37').Value = '  -2.02%  '
$ws.Range('E38').Value = '  -2.28%  '
$ws.Range('D39').Value = '347.88'
$ws.Range('E39').Value = '  +1.87%  '
$ws.Range('D40').Value = '6.34'
$ws.Range('E40').Value = '  +0.29%  '
$ws.Range('E41').Value = '  -2.60%  '
$ws.Range('E42').Value = '  -3.01%  '
$ws.Range('D43').Value = '21.58'
$ws.Range('E43').Value = '  -4.50%  '
$ws.Range('D44').Value = '21.98'
$ws.Range('E44').Value = '  -2.86%  '
$ws.Range('E45').Value = '  -3.82%  '
$ws.Range('E46').Value = '  -1.08%  '
$ws.Range('E47').Value = '  -3.58%  '
$ws.Range('D48').Value = '0.627'
$ws.Range('E48').Value = '  -4.16%  '
$ws.Range('E49').Value = '  -2.25%  '
$ws.Range('E50').Value = '  -0.01%  '
$ws.Range('E51').Value = '  +0.18%  '
